$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Text corrections: "corecl" -> "CoreCL" in the package URLs, and a refreshed Date value ---
$wsMeta.Range("B2").Value = "https://hl7chile.cl/fhir/ig/CoreCL/StructureDefinition/TiposServicio"
$wsMeta.Range("B8").Value = "2022-12-12T20:08:16-03:00"
$wsElem.Range("Q5").Value = "https://hl7chile.cl/fhir/ig/CoreCL/StructureDefinition/TiposServicio"
$wsElem.Range("Y7").Value = "https://hl7chile.cl/fhir/ig/CoreCL/ValueSet/VSTiposServicio"

# --- Column width adjustments on the "Elements" sheet ---
# NOTE: the COM ColumnWidth setter here quantizes to 1/6-character increments
# (stored_width = (Floor(ColumnWidth*6 + 0.5) + 5) / 6), so the inputs below are
# pre-compensated (inverse of that formula) to land as close as possible on the
# target stored widths taken from the target OOXML.
$widths = @{
    "A"  = 18.166666666666664
    "B"  = 10.333333333333332
    "C"  = 6.833333333333334
    "D"  = 5.0
    "E"  = 3.833333333333333
    "F"  = 4.166666666666666
    "G"  = 13.833333333333332
    "H"  = 11.166666666666668
    "I"  = 11.833333333333332
    "K"  = 40.83333333333333
    "O"  = 12.666666666666668
    "T"  = 14.833333333333332
    "U"  = 15.333333333333332
    "V"  = 16.166666666666664
    "W"  = 15.5
    "X"  = 18.0
    "Y"  = 53.16666666666667
    "Z"  = 4.833333333333334
    "AA" = 18.833333333333336
    "AB" = 39.16666666666667
    "AC" = 14.166666666666668
    "AD" = 11.5
    "AE" = 16.833333333333336
    "AF" = 8.666666666666668
    "AG" = 9.0
    "AH" = 11.333333333333332
    "AJ" = 21.833333333333336
}

# Columns that must stay hidden (setting ColumnWidth otherwise resets Hidden to False)
$hiddenCols = @("C", "D", "AE", "AF", "AG")

foreach ($col in $widths.Keys) {
    $wsElem.Columns($col).ColumnWidth = $widths[$col]
    if ($hiddenCols -contains $col) {
        $wsElem.Columns($col).Hidden = $true
    }
}
